$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 5 data (same pattern as row 4, with incremented Pull # and same cable selection)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "2C#2"
$ws.Range("C5").Value = "EXPRESS"
$ws.Range("D5").Value = "100+00"
$ws.Range("E5").Value = "200+00"

# Update the selection as in the diff (activeCell E8, sqref E8)
$ws.Range("E8").Select()
